$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Variant"
$ws.Range("B1").Value = "Stock "
